$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Localización" column (D) which held values like "41.5N35.99W"
$ws.Columns("D").Delete()

# Insert two new columns before the "Tipo" column (now column D) for Latitud/Longitud
$ws.Columns("D:E").Insert()

# Headers
$ws.Range("D1").Value = "Latitud"
$ws.Range("E1").Value = "Longitud"

# Values
$ws.Range("D2").Value = 12.36
$ws.Range("E2").Value = 15.265000000000001

# Update selection to match the target state
$ws.Range("D2").Select() | Out-Null
